# Rename "Shrub" vegetation treatment to "CSS" in the groups column.
# Final column values (A1:A9) after the edit, matching the target workbook state.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = "groups"
$ws.Range("A2").Value = "0 x Grassland"
$ws.Range("A3").Value = "6 x CSS"
$ws.Range("A4").Value = "3 x CSS"
$ws.Range("A5").Value = "5 x Grassland"
$ws.Range("A6").Value = "5 x CSS"
$ws.Range("A7").Value = "0 x CSS"
$ws.Range("A8").Value = "3 x Grassland"
$ws.Range("A9").Value = "6 x Grassland"
